$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1) -- date advances by one day; header labels unchanged.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = 45750

# ---------------------------------------------------------------------------
# 2) Bet rows. Replace the two groups of three bets (and their "CUOTA TOTAL"
#    subtotal rows) with the new set of bets / results, reusing the existing
#    row layout and visual formatting.
# ---------------------------------------------------------------------------

# --- Group 1: MIL vs PHI (rows 2-4), subtotal row 5 ---
$ws.Range("A2").Value = "MIL vs PHI"
$ws.Range("B2").Value = "Kyle Kuzma - Debajo 5.5 Dobles realizados"
$ws.Range("C2").Value = "1,42x"
$ws.Range("D2").Value = "Supero la linea en 2/10 https://prnt.sc/-3BUqD5x4j4q"

$ws.Range("A3").Value = "MIL vs PHI"
$ws.Range("B3").Value = "Giannis Antetokounmpo - Sobre 5.5 Tiros libres realizados"
$ws.Range("C3").Value = "1,38x"
$ws.Range("D3").Value = "Cumplio en 10/10 https://prnt.sc/jlLBVcNBS-x0 (Opcional: -11,5 Dobles realizados 1,50x 8/10 no supero."

$ws.Range("A4").Value = "MIL vs PHI"
$ws.Range("B4").Value = "Ricky Council IV - Sobre 4.5 Dobles intentados"
$ws.Range("C4").Value = "1,58x"
$ws.Range("D4").Value = "Cumplio en 7/10 https://prnt.sc/pSfQFYnteu3T (pueden bajar lineas)"

$ws.Range("B5").Value = "CUOTA TOTAL ="
$ws.Range("C5").Value = "3,20x"

# --- Group 2: TOR vs POR (rows 6-8), subtotal row 9 ---
$ws.Range("A6").Value = "TOR vs POR"
$ws.Range("B6").Value = "Scottie Barnes - Debajo 6.5 Dobles realizados"
$ws.Range("C6").Value = "1,43x"
$ws.Range("D6").Value = "Supero en 2/10 https://prnt.sc/JmB2tMRUUWfI (hizo 10 y 7)"

$ws.Range("A7").Value = "TOR vs POR"
$ws.Range("B7").Value = "Deni Avdija - Sobre 5.5 Triples intentados"
$ws.Range("C7").Value = "1,39x"
$ws.Range("D7").Value = "Cumplio en 10/10 https://prnt.sc/obxxUOTRkghj"

$ws.Range("A8").Value = "TOR vs POR"
$ws.Range("B8").Value = "Toumani Camara - Sobre 2.5 Dobles realizados"
$ws.Range("C8").Value = "1,52x"
$ws.Range("D8").Value = "Cumplio en 7/10 https://prnt.sc/mQabbqX7ZK5R"

$ws.Range("B9").Value = "CUOTA TOTAL ="
$ws.Range("C9").Value = "3,40x"

# --- Final combined odds row 10, with updated hyperlink ---
$ws.Range("B10").Value = "CUOTA TOTAL DUPLINHA) ="
$ws.Range("C10").Value = "10,88x"

$link = $ws.Range("D10")
$newUrl = "https://stake.com/sports/home?operation=withdraw&betId=afbc8a34-2c0e-4bf9-8b01-40a54ef4b893&modal=bet"
$link.Value = $newUrl
$link.Hyperlinks.Delete()
$ws.Hyperlinks.Add($link, $newUrl) | Out-Null

# ---------------------------------------------------------------------------
# 3) Formatting touch-ups that accompany the new content.
#    - B2 gains vertical-center alignment (keeps its plain font/no fill look).
#    - B3,B4,B6,B7,B8 switch from top-aligned wrap to vertical-center wrap,
#      matching the rest of the template rows.
#    - Row 3 grows taller to fit the longer wrapped text.
# ---------------------------------------------------------------------------
$ws.Range("B2").VerticalAlignment = -4108
foreach ($addr in @("B3", "B4", "B6", "B7", "B8")) {
    $ws.Range($addr).VerticalAlignment = -4108
}

$ws.Rows.Item(3).RowHeight = 36.75

# ---------------------------------------------------------------------------
# 4) Selection cursor moved to B13 in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("B13").Select() | Out-Null
